# "input and valid login added"
# Rename Sheet1 -> ValidLogin, clear the old sample data (cols C:D, row 3),
# and write the new username/password login values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Name = "ValidLogin"

# Write values in the same order the original author entered them so the
# shared-string table comes out in the same (Admin, admin123, Username,
# Password) order as the saved workbook.
$ws.Range("A2").Value = "Admin"
$ws.Range("B2").Value = "admin123"
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"

# Remove the old extra data (columns C:D and row 3) so the used range
# shrinks back down to A1:B2.
$ws.Range("C1:D3").ClearContents()
$ws.Range("A3:D3").ClearContents()

# Select row 7 (matches the saved selection A7:XFD7) before leaving the sheet.
$ws.Rows.Item(7).Select() | Out-Null
